{"js": "// Replace the 100 arithmetic-problem cells (20 rows x 5 cols) in the single\n// table of the worksheet with a new set of problems, in document order\n// (row-major), matching the author's commit \"Update master to output\n// generated at aa3dc9e\".\nconst newValues = [\n  [\"70-11=\", \"48+27=\", \"48+14=\", \"66-28=\", \"70-5=\"],\n  [\"83-55=\", \"74+18=\", \"56-49=\", \"81-66=\", \"64-48=\"],\n  [\"9+66=\", \"7+75=\", \"46-37=\", \"41-39=\", \"73-15=\"],\n  [\"94-65=\", \"61-13=\", \"13+78=\", \"31-13=\", \"91-39=\"],\n  [\"87+5=\", \"60-58=\", \"79+14=\", \"8+27=\", \"68+15=\"],\n  [\"65-17=\", \"18+43=\", \"7+46=\", \"19+22=\", \"5+88=\"],\n  [\"29+12=\", \"82-68=\", \"80-21=\", \"66+15=\", \"90-77=\"],\n  [\"93-24=\", \"85-68=\", \"27+37=\", \"65+27=\", \"19+49=\"],\n  [\"17+14=\", \"9+49=\", \"95-76=\", \"8+79=\", \"28+69=\"],\n  [\"47+45=\", \"41-19=\", \"30-13=\", \"20-7=\", \"76+9=\"],\n  [\"59+39=\", \"49+33=\", \"17+45=\", \"69+19=\", \"28+36=\"],\n  [\"68-49=\", \"33+19=\", \"48+34=\", \"64-37=\", \"97-48=\"],\n  [\"56+16=\", \"64-6=\", \"46+9=\", \"91-35=\", \"37+56=\"],\n  [\"86+6=\", \"82-13=\", \"96-39=\", \"42-39=\", \"54-29=\"],\n  [\"8+43=\", \"16-9=\", \"92-73=\", \"88+5=\", \"75-16=\"],\n  [\"24-17=\", \"75-38=\", \"9+77=\", \"84-75=\", \"84-48=\"],\n  [\"50-1=\", \"26+49=\", \"12+19=\", \"44-25=\", \"49+18=\"],\n  [\"49+3=\", \"81-73=\", \"83-17=\", \"17+58=\", \"82-43=\"],\n  [\"58+26=\", \"56+7=\", \"66+29=\", \"92-54=\", \"42-9=\"],\n  [\"95-39=\", \"19+54=\", \"5+68=\", \"96-69=\", \"38+34=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let r = 0; r < rows.items.length; r++) {\n  const row = rows.items[r];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  const rowValues = newValues[r];\n  if (!rowValues) continue;\n\n  for (let c = 0; c < cells.items.length; c++) {\n    if (rowValues[c] === undefined) continue;\n    cells.items[c].value = rowValues[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem cells (20 rows x 5 cols) in the single\n# table of the document with a new set of problems, in document order\n# (row-major), matching the author's commit \"Update master to output\n# generated at aa3dc9e\".\n\n$newValues = @(\n    @(\"70-11=\", \"48+27=\", \"48+14=\", \"66-28=\", \"70-5=\"),\n    @(\"83-55=\", \"74+18=\", \"56-49=\", \"81-66=\", \"64-48=\"),\n    @(\"9+66=\", \"7+75=\", \"46-37=\", \"41-39=\", \"73-15=\"),\n    @(\"94-65=\", \"61-13=\", \"13+78=\", \"31-13=\", \"91-39=\"),\n    @(\"87+5=\", \"60-58=\", \"79+14=\", \"8+27=\", \"68+15=\"),\n    @(\"65-17=\", \"18+43=\", \"7+46=\", \"19+22=\", \"5+88=\"),\n    @(\"29+12=\", \"82-68=\", \"80-21=\", \"66+15=\", \"90-77=\"),\n    @(\"93-24=\", \"85-68=\", \"27+37=\", \"65+27=\", \"19+49=\"),\n    @(\"17+14=\", \"9+49=\", \"95-76=\", \"8+79=\", \"28+69=\"),\n    @(\"47+45=\", \"41-19=\", \"30-13=\", \"20-7=\", \"76+9=\"),\n    @(\"59+39=\", \"49+33=\", \"17+45=\", \"69+19=\", \"28+36=\"),\n    @(\"68-49=\", \"33+19=\", \"48+34=\", \"64-37=\", \"97-48=\"),\n    @(\"56+16=\", \"64-6=\", \"46+9=\", \"91-35=\", \"37+56=\"),\n    @(\"86+6=\", \"82-13=\", \"96-39=\", \"42-39=\", \"54-29=\"),\n    @(\"8+43=\", \"16-9=\", \"92-73=\", \"88+5=\", \"75-16=\"),\n    @(\"24-17=\", \"75-38=\", \"9+77=\", \"84-75=\", \"84-48=\"),\n    @(\"50-1=\", \"26+49=\", \"12+19=\", \"44-25=\", \"49+18=\"),\n    @(\"49+3=\", \"81-73=\", \"83-17=\", \"17+58=\", \"82-43=\"),\n    @(\"58+26=\", \"56+7=\", \"66+29=\", \"92-54=\", \"42-9=\"),\n    @(\"95-39=\", \"19+54=\", \"5+68=\", \"96-69=\", \"38+34=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
